$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '47.373.62'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.63%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.496.67'
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.528'
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.545'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.33%  '
$ws.Range("E11").Value = '  +1.63%  '
$ws.Range("E12").Value = '  +1.33%  '
$ws.Range("E13").Value = '  +1.26%  '
$ws.Range("E14").Value = '  +1.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.886.87'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.86%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.500.71'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.854'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '47.312.31'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.94'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.36%  '
$ws.Range("E21").Value = '  +2.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.72'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '251.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.58%  '
$ws.Range("E25").Value = '  +3.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.74%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  +4.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '35.15'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.65%  '
$ws.Range("E30").Value = '  +10.39%  '
$ws.Range("E31").Value = '  -8.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.46'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.51'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.76'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.52%  '
$ws.Range("E35").Value = '  +4.14%  '
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.99'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.69'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.97%  '
$ws.Range("E40").Value = '  +2.02%  '
$ws.Range("E41").Value = '  +1.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '121.79'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.15'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.68%  '
$ws.Range("E44").Value = '  +3.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.967.69'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.58%  '
$ws.Range("E46").Value = '  +2.86%  '
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.80'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.08'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("E50").Value = '  +10.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.75'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.20%  '
